$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Rule label in B11 from "R40" to "1" (stored as text)
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
